$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values in columns A (Score for Negativity) and B (Percent Dominance)
$ws.Range("A2").Value = -0.0708726054971081
$ws.Range("B2").Value = 0.09541984732824428

$ws.Range("A3").Value = 0.1400158327575333

$ws.Range("A4").Value = -0.01495059315321198
$ws.Range("B4").Value = 0.1206106870229008

$ws.Range("A5").Value = -0.3101741051898327
$ws.Range("B5").Value = 0.05190839694656488

$ws.Range("A6").Value = -0.2200793847966849
$ws.Range("B6").Value = 0.1679389312977099

$ws.Range("A7").Value = 0.4761604127767897
$ws.Range("B7").Value = 0.2206106870229008

# Update "Terms" text column (D) to reflect re-ordered term lists
$ws.Range("D4").Value = "payment,interest,pay,amount,rate,year,month,paid,paying,balance,principal,time,made,monthly,one,making,applied,owe,make,money"
$ws.Range("D5").Value = "forbearance,told,received,deferment,form,letter,month,said,could,application,back,sent,called,repayment,time,paperwork,stating,rep,year,college"
$ws.Range("D6").Value = "account,payment,customer,credit,service,check,called,sent,received,information,told,representative,bank,letter,issue,report,back,could,day,delinquent"
$ws.Range("D7").Value = "call,phone,day,told,time,said,called,number,calling,asked,payment,even,back,know,one,never,person,got,stop,someone"
